$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 19330051920292
$ws.Range("A3").Value = 19330051920295

$ws.Range("B2").Value = "RIVERA"
$ws.Range("B3").Value = "SANCHEZ"

$ws.Range("C2").Value = "FLORES"
$ws.Range("C3").Value = "TEZOCO"

$ws.Range("D2").Value = "KARLA"
$ws.Range("D3").Value = "ESMERALDA"

$ws.Range("E2").Value = "REALIZA ANÁLISIS CITOQUÍMICOS A LÍQUIDOS Y SECRECIONES CORPORALES"
$ws.Range("E3").Value = "REALIZA ANÁLISIS CITOQUÍMICOS A LÍQUIDOS Y SECRECIONES CORPORALES"

$ws.Range("F2").Value = "4ALCV"
$ws.Range("F3").Value = "4ALCV"

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
